# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.393.23"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.239.12"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.93"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.29"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.241.61"
$ws.Range("E8").Value = "  +3.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.43"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.70"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.761.55"
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.218.98"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.397.80"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.80"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.91"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.21"
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.93"
$ws.Range("E23").Value = "  +3.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.45"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.29"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.30"
$ws.Range("E28").Value = "  +4.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.12"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("E30").Value = "  +3.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.57"
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.58"
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.91"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.98"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0732"
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "425.43"
$ws.Range("E40").Value = "  -3.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.81"
$ws.Range("E41").Value = "  -5.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.44"
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.987.55"
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("E44").Value = "  -7.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.269"
$ws.Range("E45").Value = "  +3.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.18"
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.04"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.29"
$ws.Range("E50").Value = "  -4.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.90"
$ws.Range("E51").Value = "  -0.73%  "
